$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-06-11 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-06-12 Wednesday", 2) | Out-Null
$d.Content.Find.Execute("322×2=644", $true, $false, $false, $false, $false, $true, 1, $false, "291×3=873", 2) | Out-Null
$d.Content.Find.Execute("537×2=1074", $true, $false, $false, $false, $false, $true, 1, $false, "730×8=5840", 2) | Out-Null
$d.Content.Find.Execute("518×6=3108", $true, $false, $false, $false, $false, $true, 1, $false, "931×8=7448", 2) | Out-Null
$d.Content.Find.Execute("276×9=2484", $true, $false, $false, $false, $false, $true, 1, $false, "137×7=959", 2) | Out-Null
$d.Content.Find.Execute("871×3=2613", $true, $false, $false, $false, $false, $true, 1, $false, "982×3=2946", 2) | Out-Null
$d.Content.Find.Execute("638×2=1276", $true, $false, $false, $false, $false, $true, 1, $false, "493×5=2465", 2) | Out-Null
$d.Content.Find.Execute("964×2=1928", $true, $false, $false, $false, $false, $true, 1, $false, "452×5=2260", 2) | Out-Null
$d.Content.Find.Execute("120×3=360", $true, $false, $false, $false, $false, $true, 1, $false, "665×8=5320", 2) | Out-Null
$d.Content.Find.Execute("373×8=2984", $true, $false, $false, $false, $false, $true, 1, $false, "289×2=578", 2) | Out-Null
$d.Content.Find.Execute("700×8=5600", $true, $false, $false, $false, $false, $true, 1, $false, "461×2=922", 2) | Out-Null
$d.Content.Find.Execute("460×6=2760", $true, $false, $false, $false, $false, $true, 1, $false, "472×8=3776", 2) | Out-Null
$d.Content.Find.Execute("740×9=6660", $true, $false, $false, $false, $false, $true, 1, $false, "995×6=5970", 2) | Out-Null
$d.Content.Find.Execute("538×9=4842", $true, $false, $false, $false, $false, $true, 1, $false, "894×5=4470", 2) | Out-Null
$d.Content.Find.Execute("846×7=5922", $true, $false, $false, $false, $false, $true, 1, $false, "866×3=2598", 2) | Out-Null
$d.Content.Find.Execute("237×3=711", $true, $false, $false, $false, $false, $true, 1, $false, "299×2=598", 2) | Out-Null
$d.Content.Find.Execute("316×5=1580", $true, $false, $false, $false, $false, $true, 1, $false, "178×8=1424", 2) | Out-Null
$d.Content.Find.Execute("551×5=2755", $true, $false, $false, $false, $false, $true, 1, $false, "945×3=2835", 2) | Out-Null
$d.Content.Find.Execute("236×2=472", $true, $false, $false, $false, $false, $true, 1, $false, "421×7=2947", 2) | Out-Null
$d.Content.Find.Execute("893×6=5358", $true, $false, $false, $false, $false, $true, 1, $false, "106×4=424", 2) | Out-Null
$d.Content.Find.Execute("754×8=6032", $true, $false, $false, $false, $false, $true, 1, $false, "266×2=532", 2) | Out-Null
$d.Content.Find.Execute("815×5=4075", $true, $false, $false, $false, $false, $true, 1, $false, "675×6=4050", 2) | Out-Null
$d.Content.Find.Execute("603×2=1206", $true, $false, $false, $false, $false, $true, 1, $false, "499×8=3992", 2) | Out-Null
$d.Content.Find.Execute("266×8=2128", $true, $false, $false, $false, $false, $true, 1, $false, "814×3=2442", 2) | Out-Null
$d.Content.Find.Execute("989×7=6923", $true, $false, $false, $false, $false, $true, 1, $false, "673×9=6057", 2) | Out-Null
$d.Content.Find.Execute("199×2=398", $true, $false, $false, $false, $false, $true, 1, $false, "490×3=1470", 2) | Out-Null
